# Updated cryptos list on Thu Aug 22 23:27:43 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.294.25'
$ws.Range("E2").Value = '  -0.82%  '

# Row 3
$ws.Range("D3").Value = '2.612.02'
$ws.Range("E3").Value = '  +0.30%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.36'
$ws.Range("E5").Value = '  +2.70%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.15'
$ws.Range("E6").Value = '  +0.18%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.16%  '

# Row 8
$ws.Range("E8").Value = '  -0.35%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.51'
$ws.Range("E9").Value = '  -0.04%  '

# Row 10
$ws.Range("E10").Value = '  -0.96%  '

# Row 11
$ws.Range("E11").Value = '  +2.17%  '

# Row 12
$ws.Range("E12").Value = '  +1.19%  '

# Row 13
$ws.Range("D13").Value = '3.072.56'

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.86'
$ws.Range("E14").Value = '  +5.82%  '

# Row 15
$ws.Range("D15").Value = '60.283.24'
$ws.Range("E15").Value = '  -0.85%  '

# Row 16
$ws.Range("E16").Value = '  -0.60%  '

# Row 17
$ws.Range("D17").Value = '2.614.25'
$ws.Range("E17").Value = '  -0.10%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.43'
$ws.Range("E18").Value = '  +2.27%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.64'
$ws.Range("E19").Value = '  -0.30%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.58'
$ws.Range("E20").Value = '  -0.74%  '

# Row 21
$ws.Range("E21").Value = '  -2.17%  '

# Row 22
$ws.Range("E22").Value = '  -0.21%  '

# Row 23
$ws.Range("E23").Value = '  +1.40%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.76'
$ws.Range("E24").Value = '  -0.86%  '

# Row 25
$ws.Range("E25").Value = '  +0.34%  '

# Row 26
$ws.Range("E26").Value = '  -0.36%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.00'
$ws.Range("E27").Value = '  +3.67%  '

# Row 28
$ws.Range("E28").Value = '  +5.97%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0798'
$ws.Range("E29").Value = '  +0.41%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.48'
$ws.Range("E30").Value = '  +2.14%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.95'
$ws.Range("E31").Value = '  +4.17%  '

# Row 32
$ws.Range("E32").Value = '  +0.13%  '

# Row 33
$ws.Range("E33").Value = '  -0.08%  '

# Row 34
$ws.Range("E34").Value = '  +6.00%  '

# Row 35
$ws.Range("E35").Value = '  +8.61%  '

# Row 36
$ws.Range("E36").Value = '  +0.85%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.64'
$ws.Range("E37").Value = '  +2.37%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '319.87'
$ws.Range("E38").Value = '  +7.43%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.44'
$ws.Range("E39").Value = '  +1.72%  '

# Row 40
$ws.Range("E40").Value = '  +3.62%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.851'
$ws.Range("E41").Value = '  +0.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '135.52'
$ws.Range("E42").Value = '  -3.39%  '

# Row 43
$ws.Range("E43").Value = '  +0.48%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.32%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.97'
$ws.Range("E45").Value = '  +2.16%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.609'
$ws.Range("E46").Value = '  +0.48%  '

# Row 47
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0552'
$ws.Range("E47").Value = '  +0.61%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.01'
$ws.Range("E48").Value = '  +2.33%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.04'
$ws.Range("E49").Value = '  +1.53%  '

# Row 50
$ws.Range("E50").Value = '  -0.32%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.76'
$ws.Range("E51").Value = '  +0.62%  '
